$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Formula = "=MEDIAN(B1:B18)"
$ws.Range("E18").Value = 0
$ws.Range("F18").Formula = "=_xlfn.STDEV.P(B1:B18)"
$ws.Range("G18").Formula = "=_xlfn.VAR.P(B1:B18)"
$ws.Range("E18").Select()
